$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial number (days since 1899-12-30).
# Every data row (2 through 451) has its value bumped from
# 45177 (2023-09-08) to 45178 (2023-09-09).
for ($row = 2; $row -le 451; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45177) {
        $cell.Value2 = 45178
    }
}
